$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and the Cosmos / WrappedliquidstakedEther2.0
# row swap) per the Thu May 25 15:23:33 UTC 2023 GitHub Actions refresh.
# A leading apostrophe forces each value to be stored as literal text, matching the
# original inline-string cells (prevents Excel from re-interpreting numeric-looking
# strings such as "1.007" or "10.51" as actual numbers).

$ws.Range("D2").Value = "'26.431.56"
$ws.Range("E2").Value = "'  -0.05%  "

$ws.Range("D3").Value = "'1.804.22"
$ws.Range("E3").Value = "'  +0.02%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "'  -0.02%  "

$ws.Range("D5").Value = "'1.006"
$ws.Range("E5").Value = "'  -0.17%  "

$ws.Range("D6").Value = "'306.25"
$ws.Range("E6").Value = "'  -0.46%  "

$ws.Range("D7").Value = "'0.4524"
$ws.Range("E7").Value = "'  -0.31%  "

$ws.Range("D8").Value = "'0.3601"
$ws.Range("E8").Value = "'  -1.19%  "

$ws.Range("D9").Value = "'46.35"
$ws.Range("E9").Value = "'  +1.95%  "

$ws.Range("D10").Value = "'0.07068"
$ws.Range("E10").Value = "'  -0.39%  "

$ws.Range("D11").Value = "'0.8925"
$ws.Range("E11").Value = "'  +2.57%  "

$ws.Range("D12").Value = "'0.07817"
$ws.Range("E12").Value = "'  +0.47%  "

$ws.Range("D13").Value = "'19.41"
$ws.Range("E13").Value = "'  +1.03%  "

$ws.Range("D14").Value = "'1.871.83"
$ws.Range("E14").Value = "'  +2.27%  "

$ws.Range("D15").Value = "'5.288"
$ws.Range("E15").Value = "'  +0.50%  "

$ws.Range("D16").Value = "'6.316"
$ws.Range("E16").Value = "'  -0.11%  "

$ws.Range("D17").Value = "'85.65"
$ws.Range("E17").Value = "'  -0.84%  "

$ws.Range("D18").Value = "'1.008"
$ws.Range("E18").Value = "'  -0.01%  "

$ws.Range("D19").Value = "'0.000008483"
$ws.Range("E19").Value = "'  -0.88%  "

$ws.Range("E20").Value = "'  -0.45%  "

$ws.Range("D21").Value = "'26.471.42"
$ws.Range("E21").Value = "'  +0.00%  "

$ws.Range("D22").Value = "'14.21"
$ws.Range("E22").Value = "'  +0.12%  "

$ws.Range("D23").Value = "'4.967"
$ws.Range("E23").Value = "'  +0.44%  "

$ws.Range("B24").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "'2.041.93"
$ws.Range("E24").Value = "'  -1.42%  "

$ws.Range("B25").Value = "'Cosmos"
$ws.Range("C25").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'10.51"
$ws.Range("E25").Value = "'  +1.73%  "

$ws.Range("D26").Value = "'1.959"
$ws.Range("E26").Value = "'  -0.63%  "

$ws.Range("D27").Value = "'151.66"
$ws.Range("E27").Value = "'  +0.77%  "

$ws.Range("D28").Value = "'17.79"
$ws.Range("E28").Value = "'  -0.24%  "

$ws.Range("D29").Value = "'2.068"
$ws.Range("E29").Value = "'  +3.64%  "

$ws.Range("D30").Value = "'112.08"
$ws.Range("E30").Value = "'  -0.84%  "

$ws.Range("D31").Value = "'4.852"
$ws.Range("E31").Value = "'  -0.14%  "

$ws.Range("D32").Value = "'0.08692"
$ws.Range("E32").Value = "'  +0.26%  "

$ws.Range("D33").Value = "'3.114"
$ws.Range("E33").Value = "'  +0.06%  "

$ws.Range("D34").Value = "'2.831"
$ws.Range("E34").Value = "'  +13.43%  "

$ws.Range("D36").Value = "'0.7235"
$ws.Range("E36").Value = "'  -0.37%  "

$ws.Range("D37").Value = "'1.106"
$ws.Range("E37").Value = "'  -0.21%  "

$ws.Range("D38").Value = "'1.075"
$ws.Range("E38").Value = "'  -0.11%  "

$ws.Range("D39").Value = "'0.01934"
$ws.Range("E39").Value = "'  +1.55%  "

$ws.Range("D40").Value = "'0.05107"
$ws.Range("E40").Value = "'  +0.49%  "

$ws.Range("D41").Value = "'2.902"
$ws.Range("E41").Value = "'  +1.64%  "

$ws.Range("D42").Value = "'0.5109"
$ws.Range("E42").Value = "'  +4.36%  "

$ws.Range("D43").Value = "'6.764"
$ws.Range("E43").Value = "'  -1.61%  "

$ws.Range("D44").Value = "'0.1514"
$ws.Range("E44").Value = "'  -3.19%  "

$ws.Range("D45").Value = "'8.019"
$ws.Range("E45").Value = "'  -1.18%  "

$ws.Range("D46").Value = "'0.4671"
$ws.Range("E46").Value = "'  +1.89%  "

$ws.Range("E47").Value = "'  -0.30%  "

$ws.Range("D48").Value = "'10.03"
$ws.Range("E48").Value = "'  +0.84%  "

$ws.Range("D49").Value = "'100.28"
$ws.Range("E49").Value = "'  -1.21%  "

$ws.Range("D50").Value = "'1.575"
$ws.Range("E50").Value = "'  +0.02%  "

$ws.Range("E51").Value = "'  -0.03%  "

